# Add model/view initial implementation for dccsettings
# Updates the "Translation" sheet (Table8) rows 5-18 with new Text IDs,
# Typography names, Alignment and GB (translation) text, and appends the
# new rows 13-18 that were previously blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 5: wildcard text id, left aligned, value placeholder
$ws.Range("B5").Value = "wildcardTextId"
$ws.Range("C5").Value = "Typography_00"
$ws.Range("D5").Value = "Left"
$ws.Range("E5").Value = "<value>"
$ws.Range("F5").Value = "LTR"

# Row 6: TOFF
$ws.Range("B6").Value = "SingleUseId14"
$ws.Range("C6").Value = "Typography_00"
$ws.Range("D6").Value = "Left"
$ws.Range("E6").Value = "TOFF"
$ws.Range("F6").Value = "LTR"

# Row 7: wildcard text id, center aligned, value placeholder
$ws.Range("B7").Value = "wildcardTextIdCenter"
$ws.Range("C7").Value = "Typography_00"
$ws.Range("D7").Value = "Center"
$ws.Range("E7").Value = "<value>"
$ws.Range("F7").Value = "LTR"

# Row 8: Trip Current
$ws.Range("B8").Value = "SingleUseId15"
$ws.Range("C8").Value = "Typography_00"
$ws.Range("D8").Value = "Left"
$ws.Range("E8").Value = "Trip Current"
$ws.Range("F8").Value = "LTR"

# Row 9: Slew Rate
$ws.Range("B9").Value = "SingleUseId16"
$ws.Range("C9").Value = "Typography_00"
$ws.Range("D9").Value = "Left"
$ws.Range("E9").Value = "Slew Rate"
$ws.Range("F9").Value = "LTR"

# Row 10: Status
$ws.Range("B10").Value = "SingleUseId17"
$ws.Range("C10").Value = "Typography_00"
$ws.Range("D10").Value = "Left"
$ws.Range("E10").Value = "Status"
$ws.Range("F10").Value = "LTR"

# Row 11: Prog.
$ws.Range("B11").Value = "SingleUseId18"
$ws.Range("C11").Value = "Typography_00"
$ws.Range("D11").Value = "Left"
$ws.Range("E11").Value = "Prog."
$ws.Range("F11").Value = "LTR"

# Row 12: Track
$ws.Range("B12").Value = "SingleUseId19"
$ws.Range("C12").Value = "Typography_00"
$ws.Range("D12").Value = "Left"
$ws.Range("E12").Value = "Track"
$ws.Range("F12").Value = "LTR"

# Row 13: Fault (new row)
$ws.Range("B13").Value = "SingleUseId20"
$ws.Range("C13").Value = "Typography_00"
$ws.Range("D13").Value = "Left"
$ws.Range("E13").Value = "Fault"
$ws.Range("F13").Value = "LTR"

# Row 14: Over Temp (new row)
$ws.Range("B14").Value = "SingleUseId21"
$ws.Range("C14").Value = "Typography_00"
$ws.Range("D14").Value = "Left"
$ws.Range("E14").Value = "Over Temp"
$ws.Range("F14").Value = "LTR"

# Row 15: Over Current (new row)
$ws.Range("B15").Value = "SingleUseId22"
$ws.Range("C15").Value = "Typography_00"
$ws.Range("D15").Value = "Left"
$ws.Range("E15").Value = "Over Current"
$ws.Range("F15").Value = "LTR"

# Row 16: Open Load (new row)
$ws.Range("B16").Value = "SingleUseId23"
$ws.Range("C16").Value = "Typography_00"
$ws.Range("D16").Value = "Left"
$ws.Range("E16").Value = "Open Load"
$ws.Range("F16").Value = "LTR"

# Row 17: Current (new row)
$ws.Range("B17").Value = "SingleUseId28"
$ws.Range("C17").Value = "Typography_00"
$ws.Range("D17").Value = "Left"
$ws.Range("E17").Value = "Current"
$ws.Range("F17").Value = "LTR"

# Row 18: Backlight (new row)
$ws.Range("B18").Value = "SingleUseId30"
$ws.Range("C18").Value = "Typography_00"
$ws.Range("D18").Value = "Left"
$ws.Range("E18").Value = "Backlight"
$ws.Range("F18").Value = "LTR"

# The newly written rows 13-18 pick up the column's default style (style
# index 1) since they previously had no cells at all. The other rows in
# this table use the workbook default ("Normal") style, so normalize the
# new cells to match and keep formatting consistent across the table.
$ws.Range("B13:F18").Style = "Normal"
